$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 210.77777
$ws.Range("I6").Value = 210.77777
$ws.Range("K6").Value = 632.33331
$ws.Range("M6").Value = -520.33331
$ws.Range("H9").Value = 1299099.8
$ws.Range("I9").Value = 3246876.8
$ws.Range("J9").Value = 581.6667
$ws.Range("K9").Value = 3246876.8
$ws.Range("L9").Value = 581.6667
$ws.Range("M9").Value = -3246707.8
$ws.Range("N9").Value = -919.6667
$ws.Range("H18").Value = 984.6667
$ws.Range("I18").Value = 984.6667
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 984.6667
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = -700.6667
$ws.Range("H33").Value = 723.90625
$ws.Range("I33").Value = 282.17242
$ws.Range("K33").Value = 282.17242
$ws.Range("M33").Value = -53.17241999999999
$ws.Range("H43").Value = 13188.5
$ws.Range("I43").Value = 16231.667
$ws.Range("K43").Value = 16231.667
$ws.Range("M43").Value = -16162.667
$ws.Range("H62").Value = 1795.6666
$ws.Range("I62").Value = 1793.5
$ws.Range("J62").Value = 1800
$ws.Range("K62").Value = 1793.5
$ws.Range("L62").Value = 1800
$ws.Range("M62").Value = -1169.5
$ws.Range("N62").Value = -3048
$ws.Range("H64").Value = 7899.5
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("H65").Value = 1795.6666
$ws.Range("I65").Value = 1793.5
$ws.Range("J65").Value = 1800
$ws.Range("K65").Value = 8967.5
$ws.Range("L65").Value = 9000
$ws.Range("M65").Value = -5847.5
$ws.Range("N65").Value = -15240
$ws.Range("H67").Value = 7899.5
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("H76").Value = 9671.190000000001
$ws.Range("I76").Value = 9348.3125
$ws.Range("J76").Value = 10704.4
$ws.Range("K76").Value = 9348.3125
$ws.Range("L76").Value = 10704.4
$ws.Range("M76").Value = -9033.3125
$ws.Range("N76").Value = -11334.4
$ws.Range("H79").Value = 9671.190000000001
$ws.Range("I79").Value = 9348.3125
$ws.Range("J79").Value = 10704.4
$ws.Range("K79").Value = 9348.3125
$ws.Range("L79").Value = 10704.4
$ws.Range("M79").Value = -8256.3125
$ws.Range("N79").Value = -12888.4
$ws.Range("H80").Value = 85086.086
$ws.Range("I80").Value = 745.8
$ws.Range("J80").Value = 145329.14
$ws.Range("K80").Value = 2237.4
$ws.Range("L80").Value = 435987.42
$ws.Range("M80").Value = -1239.4
$ws.Range("N80").Value = -437983.42
$ws.Range("H83").Value = 85086.086
$ws.Range("I83").Value = 745.8
$ws.Range("J83").Value = 145329.14
$ws.Range("K83").Value = 6712.2
$ws.Range("L83").Value = 1307962.26
$ws.Range("M83").Value = -1720.2
$ws.Range("N83").Value = -1317946.26
$ws.Range("H97").Value = 2676.3333
$ws.Range("J97").Value = 2676.3333
$ws.Range("L97").Value = 8028.999899999999
$ws.Range("N97").Value = -9020.999899999999
$ws.Range("H107").Value = 1351.0667
$ws.Range("I107").Value = 1243.5
$ws.Range("J107").Value = 1566.2
$ws.Range("K107").Value = 1243.5
$ws.Range("L107").Value = 1566.2
$ws.Range("M107").Value = 676.5
$ws.Range("N107").Value = -5406.2
$ws.Range("H111").Value = 10150.536
$ws.Range("I111").Value = 2324.4
$ws.Range("K111").Value = 6973.200000000001
$ws.Range("M111").Value = -3906.200000000001
$ws.Range("H112").Value = 5641.8
$ws.Range("J112").Value = 8503
$ws.Range("L112").Value = 25509
$ws.Range("N112").Value = -27725
$ws.Range("H116").Value = 6983.4707
$ws.Range("I116").Value = 6504.5
$ws.Range("J116").Value = 7409.222
$ws.Range("K116").Value = 6504.5
$ws.Range("L116").Value = 7409.222
$ws.Range("M116").Value = -3062.5
$ws.Range("N116").Value = -14293.222
$ws.Range("H132").Value = 3449.44
$ws.Range("I132").Value = 3625.6191
$ws.Range("J132").Value = 2524.5
$ws.Range("K132").Value = 10876.8573
$ws.Range("L132").Value = 7573.5
$ws.Range("M132").Value = -8346.8573
$ws.Range("N132").Value = -12633.5
$ws.Range("H136").Value = 87994.5
$ws.Range("J136").Value = 87994.5
$ws.Range("L136").Value = 87994.5
$ws.Range("N136").Value = -98194.5
$ws.Range("H141").Value = 3758.4614
$ws.Range("I141").Value = 3287.1
$ws.Range("J141").Value = 5329.6665
$ws.Range("K141").Value = 9861.299999999999
$ws.Range("L141").Value = 15988.9995
$ws.Range("M141").Value = -4681.299999999999
$ws.Range("N141").Value = -26348.9995
$ws.Range("N18").ClearContents()
$ws.Range("N64").ClearContents()
$ws.Range("N67").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1713.6923
$ws.Range("I2").Value = 2209.2222
$ws.Range("J2").Value = 598.75
$ws.Range("K2").Value = 2209.2222
$ws.Range("L2").Value = 598.75
$ws.Range("M2").Value = -2096.2222
$ws.Range("N2").Value = -824.75
$ws.Range("H4").Value = 1039.6
$ws.Range("J4").Value = 1714.4
$ws.Range("L4").Value = 1714.4
$ws.Range("N4").Value = -1946.4
$ws.Range("H12").Value = 300
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 300
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 300
$ws.Range("N12").Value = -646
$ws.Range("H32").Value = 1907.9861
$ws.Range("I32").Value = 788.07245
$ws.Range("K32").Value = 788.07245
$ws.Range("M32").Value = -501.07245
$ws.Range("H33").Value = 48891.668
$ws.Range("I33").Value = 48891.668
$ws.Range("K33").Value = 48891.668
$ws.Range("M33").Value = -48562.668
$ws.Range("H45").Value = 12603.315
$ws.Range("I45").Value = 13870.823
$ws.Range("K45").Value = 13870.823
$ws.Range("M45").Value = -13493.823
$ws.Range("H46").Value = 9177.833000000001
$ws.Range("I46").Value = 6871.5
$ws.Range("J46").Value = 10331
$ws.Range("K46").Value = 6871.5
$ws.Range("L46").Value = 10331
$ws.Range("M46").Value = -6552.5
$ws.Range("N46").Value = -10969
$ws.Range("H61").Value = 6824.706
$ws.Range("I61").Value = 6681.0835
$ws.Range("J61").Value = 7169.4
$ws.Range("K61").Value = 6681.0835
$ws.Range("L61").Value = 7169.4
$ws.Range("M61").Value = -6469.0835
$ws.Range("N61").Value = -7593.4
$ws.Range("H74").Value = 1467.6666
$ws.Range("I74").Value = 1391
$ws.Range("J74").Value = 1851
$ws.Range("K74").Value = 1391
$ws.Range("L74").Value = 1851
$ws.Range("M74").Value = -517
$ws.Range("N74").Value = -3599
$ws.Range("H77").Value = 1467.6666
$ws.Range("I77").Value = 1391
$ws.Range("J77").Value = 1851
$ws.Range("K77").Value = 6955
$ws.Range("L77").Value = 9255
$ws.Range("M77").Value = -2587
$ws.Range("N77").Value = -17991
$ws.Range("H92").Value = 60000
$ws.Range("J92").Value = 60000
$ws.Range("L92").Value = 60000
$ws.Range("N92").Value = -64992
$ws.Range("H102").Value = 4477.25
$ws.Range("I102").Value = 4175
$ws.Range("K102").Value = 4175
$ws.Range("M102").Value = -2553
$ws.Range("H110").Value = 2643.0715
$ws.Range("I110").Value = 2131.8
$ws.Range("J110").Value = 3921.25
$ws.Range("K110").Value = 2131.8
$ws.Range("L110").Value = 3921.25
$ws.Range("M110").Value = -86.80000000000018
$ws.Range("N110").Value = -8011.25
$ws.Range("H116").Value = 1713.6923
$ws.Range("I116").Value = 2209.2222
$ws.Range("J116").Value = 598.75
$ws.Range("K116").Value = 2209.2222
$ws.Range("L116").Value = 598.75
$ws.Range("M116").Value = 84.77779999999984
$ws.Range("N116").Value = -5186.75
$ws.Range("H125").Value = 82374
$ws.Range("J125").Value = 82374
$ws.Range("L125").Value = 82374
$ws.Range("N125").Value = -92214
$ws.Range("H132").Value = 2263.7073
$ws.Range("I132").Value = 2046.1794
$ws.Range("J132").Value = 6505.5
$ws.Range("K132").Value = 6138.5382
$ws.Range("L132").Value = 19516.5
$ws.Range("M132").Value = -3608.5382
$ws.Range("N132").Value = -24576.5
$ws.Range("H136").Value = 6824.706
$ws.Range("I136").Value = 6681.0835
$ws.Range("J136").Value = 7169.4
$ws.Range("K136").Value = 20043.2505
$ws.Range("L136").Value = 21508.2
$ws.Range("M136").Value = -17493.2505
$ws.Range("N136").Value = -26608.2
$ws.Range("M12").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1713.6923
$ws.Range("I3").Value = 2209.2222
$ws.Range("J3").Value = 598.75
$ws.Range("K3").Value = 2209.2222
$ws.Range("L3").Value = 598.75
$ws.Range("M3").Value = -2095.2222
$ws.Range("N3").Value = -826.75
$ws.Range("H16").Value = 1899
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 1899
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 1899
$ws.Range("N16").Value = -2239
$ws.Range("H20").Value = 2062.44
$ws.Range("I20").Value = 2382.2666
$ws.Range("J20").Value = 1582.7
$ws.Range("K20").Value = 2382.2666
$ws.Range("L20").Value = 1582.7
$ws.Range("M20").Value = -2135.2666
$ws.Range("N20").Value = -2076.7
$ws.Range("H64").Value = 3086.4443
$ws.Range("J64").Value = 3118
$ws.Range("L64").Value = 3118
$ws.Range("N64").Value = -3568
$ws.Range("H67").Value = 3086.4443
$ws.Range("J67").Value = 3118
$ws.Range("L67").Value = 3118
$ws.Range("N67").Value = -4678
$ws.Range("H88").Value = 35497.8
$ws.Range("I88").Value = 15244.5
$ws.Range("J88").Value = 49000
$ws.Range("K88").Value = 15244.5
$ws.Range("L88").Value = 49000
$ws.Range("M88").Value = -14838.5
$ws.Range("N88").Value = -49812
$ws.Range("H91").Value = 35497.8
$ws.Range("I91").Value = 15244.5
$ws.Range("J91").Value = 49000
$ws.Range("K91").Value = 15244.5
$ws.Range("L91").Value = 49000
$ws.Range("M91").Value = -13840.5
$ws.Range("N91").Value = -51808
$ws.Range("H99").Value = 2608.077
$ws.Range("I99").Value = 2490.55
$ws.Range("J99").Value = 2999.8333
$ws.Range("K99").Value = 2490.55
$ws.Range("L99").Value = 2999.8333
$ws.Range("M99").Value = -992.5500000000002
$ws.Range("N99").Value = -5995.8333
$ws.Range("H134").Value = 2304.7585
$ws.Range("I134").Value = 2378.423
$ws.Range("K134").Value = 7135.268999999999
$ws.Range("M134").Value = -4600.268999999999
$ws.Range("M16").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 216.66667
$ws.Range("I7").Value = 100
$ws.Range("J7").Value = 450
$ws.Range("K7").Value = 100
$ws.Range("L7").Value = 450
$ws.Range("M7").Value = 13
$ws.Range("N7").Value = -676
$ws.Range("H16").Value = 61564.383
$ws.Range("I16").Value = 22299.428
$ws.Range("J16").Value = 107373.5
$ws.Range("K16").Value = 22299.428
$ws.Range("L16").Value = 107373.5
$ws.Range("M16").Value = -22012.428
$ws.Range("N16").Value = -107947.5
$ws.Range("H22").Value = 1279.0834
$ws.Range("I22").Value = 1279.0834
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 1279.0834
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -929.0834
$ws.Range("H31").Value = 3949.0256
$ws.Range("I31").Value = 3937.037
$ws.Range("J31").Value = 3976
$ws.Range("K31").Value = 3937.037
$ws.Range("L31").Value = 3976
$ws.Range("M31").Value = -3642.037
$ws.Range("N31").Value = -4566
$ws.Range("H34").Value = 3949.0256
$ws.Range("I34").Value = 3937.037
$ws.Range("J34").Value = 3976
$ws.Range("K34").Value = 3937.037
$ws.Range("L34").Value = 3976
$ws.Range("M34").Value = -3735.037
$ws.Range("N34").Value = -4380
$ws.Range("H51").Value = 16722
$ws.Range("J51").Value = 14583.167
$ws.Range("L51").Value = 14583.167
$ws.Range("N51").Value = -16055.167
$ws.Range("H58").Value = 2513.7144
$ws.Range("I58").Value = 2534.3333
$ws.Range("K58").Value = 2534.3333
$ws.Range("M58").Value = -2331.3333
$ws.Range("H60").Value = 8888.888999999999
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 8888.888999999999
$ws.Range("K60").Value = 0
$ws.Range("L60").Value = 8888.888999999999
$ws.Range("N60").Value = -9910.888999999999
$ws.Range("H61").Value = 16722
$ws.Range("J61").Value = 14583.167
$ws.Range("L61").Value = 14583.167
$ws.Range("N61").Value = -15279.167
$ws.Range("H62").Value = 2611.5
$ws.Range("I62").Value = 2649
$ws.Range("K62").Value = 2649
$ws.Range("M62").Value = -2025
$ws.Range("H65").Value = 2611.5
$ws.Range("I65").Value = 2649
$ws.Range("K65").Value = 13245
$ws.Range("M65").Value = -10125
$ws.Range("H99").Value = 6419.8
$ws.Range("I99").Value = 5051
$ws.Range("K99").Value = 5051
$ws.Range("M99").Value = -3553
$ws.Range("H107").Value = 754.1923
$ws.Range("I107").Value = 641.6818
$ws.Range("J107").Value = 1373
$ws.Range("K107").Value = 641.6818
$ws.Range("L107").Value = 1373
$ws.Range("M107").Value = 1278.3182
$ws.Range("N107").Value = -5213
$ws.Range("H113").Value = 61564.383
$ws.Range("I113").Value = 22299.428
$ws.Range("J113").Value = 107373.5
$ws.Range("K113").Value = 22299.428
$ws.Range("L113").Value = 107373.5
$ws.Range("M113").Value = -20129.428
$ws.Range("N113").Value = -111713.5
$ws.Range("H126").Value = 6419.8
$ws.Range("I126").Value = 5051
$ws.Range("K126").Value = 15153
$ws.Range("M126").Value = -12683
$ws.Range("H132").Value = 3663.8538
$ws.Range("I132").Value = 3389.3242
$ws.Range("J132").Value = 6203.25
$ws.Range("K132").Value = 10167.9726
$ws.Range("L132").Value = 18609.75
$ws.Range("M132").Value = -7637.972600000001
$ws.Range("N132").Value = -23669.75
$ws.Range("H134").Value = 5445.7954
$ws.Range("I134").Value = 4832.027
$ws.Range("J134").Value = 8690
$ws.Range("K134").Value = 14496.081
$ws.Range("L134").Value = 26070
$ws.Range("M134").Value = -11961.081
$ws.Range("N134").Value = -31140
$ws.Range("H136").Value = 2513.7144
$ws.Range("I136").Value = 2534.3333
$ws.Range("K136").Value = 7602.999899999999
$ws.Range("M136").Value = -5052.999899999999
$ws.Range("N22").ClearContents()
$ws.Range("M60").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 1367.6
$ws.Range("J2").Value = 560
$ws.Range("L2").Value = 3360
$ws.Range("N2").Value = -3586
$ws.Range("H5").Value = 688.2
$ws.Range("I5").Value = 610.25
$ws.Range("K5").Value = 1830.75
$ws.Range("M5").Value = -1718.75
$ws.Range("H9").Value = 9332.833000000001
$ws.Range("J9").Value = 9999.25
$ws.Range("L9").Value = 29997.75
$ws.Range("N9").Value = -30445.75
$ws.Range("H34").Value = 2875.2917
$ws.Range("I34").Value = 202
$ws.Range("J34").Value = 2991.5217
$ws.Range("K34").Value = 606
$ws.Range("L34").Value = 8974.5651
$ws.Range("M34").Value = -522
$ws.Range("N34").Value = -9142.5651
$ws.Range("H38").Value = 199.5
$ws.Range("I38").Value = 199
$ws.Range("J38").Value = 200
$ws.Range("K38").Value = 597
$ws.Range("L38").Value = 600
$ws.Range("M38").Value = -250
$ws.Range("N38").Value = -1294
$ws.Range("H75").Value = 2522.5
$ws.Range("J75").Value = 2606.6667
$ws.Range("L75").Value = 7820.000100000001
$ws.Range("N75").Value = -9816.000100000001
$ws.Range("H78").Value = 2522.5
$ws.Range("J78").Value = 2606.6667
$ws.Range("L78").Value = 23460.0003
$ws.Range("N78").Value = -33444.0003
$ws.Range("H80").Value = 3994.25
$ws.Range("J80").Value = 3987.5
$ws.Range("L80").Value = 11962.5
$ws.Range("N80").Value = -13834.5
$ws.Range("H83").Value = 3994.25
$ws.Range("J83").Value = 3987.5
$ws.Range("L83").Value = 35887.5
$ws.Range("N83").Value = -45247.5
$ws.Range("H103").Value = 1099.75
$ws.Range("I103").Value = 2000
$ws.Range("J103").Value = 799.6667
$ws.Range("K103").Value = 6000
$ws.Range("L103").Value = 2399.0001
$ws.Range("M103").Value = -5121
$ws.Range("N103").Value = -4157.0001
$ws.Range("H135").Value = 688.2
$ws.Range("I135").Value = 610.25
$ws.Range("K135").Value = 5492.25
$ws.Range("M135").Value = -2957.25
$ws.Range("H141").Value = 5082.5
$ws.Range("I141").Value = 2624.25
$ws.Range("K141").Value = 7872.75
$ws.Range("M141").Value = -2692.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H9").Value = 2000
$ws.Range("J9").Value = 2000
$ws.Range("L9").Value = 2000
$ws.Range("N9").Value = -2340
$ws.Range("H80").Value = 4958.8
$ws.Range("I80").Value = 1698.6666
$ws.Range("J80").Value = 9849
$ws.Range("K80").Value = 1698.6666
$ws.Range("L80").Value = 9849
$ws.Range("M80").Value = -700.6666
$ws.Range("N80").Value = -11845
$ws.Range("H83").Value = 4958.8
$ws.Range("I83").Value = 1698.6666
$ws.Range("J83").Value = 9849
$ws.Range("K83").Value = 8493.333000000001
$ws.Range("L83").Value = 49245
$ws.Range("M83").Value = -3501.333000000001
$ws.Range("N83").Value = -59229
$ws.Range("H108").Value = 99998
$ws.Range("I108").Value = 99998
$ws.Range("J108").Value = 0
$ws.Range("K108").Value = 99998
$ws.Range("L108").Value = 0
$ws.Range("M108").Value = -96158
$ws.Range("H113").Value = 1908.6154
$ws.Range("I113").Value = 1763.8334
$ws.Range("J113").Value = 2032.7142
$ws.Range("K113").Value = 1763.8334
$ws.Range("L113").Value = 2032.7142
$ws.Range("M113").Value = 406.1666
$ws.Range("N113").Value = -6372.7142
$ws.Range("H126").Value = 4107.3438
$ws.Range("I126").Value = 3699.95
$ws.Range("J126").Value = 4786.3335
$ws.Range("K126").Value = 11099.85
$ws.Range("L126").Value = 14359.0005
$ws.Range("M126").Value = -8629.849999999999
$ws.Range("N126").Value = -19299.0005
$ws.Range("H132").Value = 2447.9556
$ws.Range("I132").Value = 2308.9268
$ws.Range("K132").Value = 6926.780400000001
$ws.Range("M132").Value = -4396.780400000001
$ws.Range("N108").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 10005000
$ws.Range("I4").Value = 10005000
$ws.Range("K4").Value = 10005000
$ws.Range("M4").Value = -10004887
$ws.Range("I16").Value = 1081.3334
$ws.Range("J16").Value = 2612.5
$ws.Range("K16").Value = 1081.3334
$ws.Range("L16").Value = 2612.5
$ws.Range("M16").Value = -911.3334
$ws.Range("N16").Value = -2952.5
$ws.Range("H19").Value = 7502.5
$ws.Range("I19").Value = 7502.5
$ws.Range("K19").Value = 7502.5
$ws.Range("M19").Value = -7332.5
$ws.Range("H22").Value = 3067.4827
$ws.Range("J22").Value = 3068.077
$ws.Range("L22").Value = 3068.077
$ws.Range("N22").Value = -3658.077
$ws.Range("H26").Value = 18004.5
$ws.Range("I26").Value = 18004.5
$ws.Range("K26").Value = 18004.5
$ws.Range("M26").Value = -17709.5
$ws.Range("H27").Value = 3067.4827
$ws.Range("J27").Value = 3068.077
$ws.Range("L27").Value = 3068.077
$ws.Range("N27").Value = -3282.077
$ws.Range("H28").Value = 10005000
$ws.Range("I28").Value = 10005000
$ws.Range("K28").Value = 10005000
$ws.Range("M28").Value = -10004768
$ws.Range("H34").Value = 10000
$ws.Range("I34").Value = 7000
$ws.Range("J34").Value = 11500
$ws.Range("K34").Value = 7000
$ws.Range("L34").Value = 11500
$ws.Range("M34").Value = -6828
$ws.Range("N34").Value = -11844
$ws.Range("H37").Value = 10005000
$ws.Range("I37").Value = 10005000
$ws.Range("K37").Value = 10005000
$ws.Range("M37").Value = -10004893
$ws.Range("H94").Value = 50000
$ws.Range("J94").Value = 50000
$ws.Range("L94").Value = 50000
$ws.Range("N94").Value = -51352
$ws.Range("H100").Value = 1791.4
$ws.Range("I100").Value = 1791.4
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 1791.4
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -1250.4
$ws.Range("H132").Value = 4288.722
$ws.Range("I132").Value = 3260.2
$ws.Range("J132").Value = 5574.375
$ws.Range("K132").Value = 9780.599999999999
$ws.Range("L132").Value = 16723.125
$ws.Range("M132").Value = -7250.599999999999
$ws.Range("N132").Value = -21783.125
$ws.Range("H136").Value = 1749.25
$ws.Range("I136").Value = 1749.25
$ws.Range("K136").Value = 5247.75
$ws.Range("M136").Value = -2697.75
$ws.Range("N100").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3116.111
$ws.Range("I81").Value = 3567.6
$ws.Range("J81").Value = 2551.75
$ws.Range("K81").Value = 7135.2
$ws.Range("L81").Value = 5103.5
$ws.Range("M81").Value = -6074.2
$ws.Range("N81").Value = -7225.5
$ws.Range("H84").Value = 3116.111
$ws.Range("I84").Value = 3567.6
$ws.Range("J84").Value = 2551.75
$ws.Range("K84").Value = 35676
$ws.Range("L84").Value = 25517.5
$ws.Range("M84").Value = -30372
$ws.Range("N84").Value = -36125.5
$ws.Range("H100").Value = 1740.4375
$ws.Range("I100").Value = 1737.3
$ws.Range("J100").Value = 1745.6666
$ws.Range("K100").Value = 3474.6
$ws.Range("L100").Value = 3491.3332
$ws.Range("M100").Value = -2933.6
$ws.Range("N100").Value = -4573.3332
$ws.Range("H113").Value = 1547.25
$ws.Range("J113").Value = 1094.5
$ws.Range("L113").Value = 3283.5
$ws.Range("N113").Value = -7623.5
$ws.Range("H118").Value = 50000
$ws.Range("J118").Value = 50000
$ws.Range("L118").Value = 50000
$ws.Range("N118").Value = -53314
$ws.Range("H132").Value = 2179.3333
$ws.Range("I132").Value = 2267.9375
$ws.Range("J132").Value = 1961.2307
$ws.Range("K132").Value = 6803.8125
$ws.Range("L132").Value = 5883.6921
$ws.Range("M132").Value = -4273.8125
$ws.Range("N132").Value = -10943.6921
$ws.Range("H136").Value = 2040.375
$ws.Range("I136").Value = 1594.2258
$ws.Range("K136").Value = 4782.6774
$ws.Range("M136").Value = -2232.6774
